# update : added CRUD Package
# Adds a new employee record ("sameer"/"sajid") to row 3 of the "empdata"
# sheet, wires up a mailto hyperlink on the email cell, stores the mobile
# number as text (quote-prefixed) and moves the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("empdata")

# New record values (fname, lname, email, mobile)
$ws.Range("A3").Value = "sameer"
$ws.Range("B3").Value = "sajid"
$ws.Range("C3").Value = "sameer@gmail.com"

# Turn the e-mail cell into a mailto: hyperlink, mirroring the existing
# C2 hyperlink already on the sheet.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:sameer@gmail.com")

# Mobile number entered as text (leading apostrophe keeps the leading
# digits/zeros intact instead of Excel coercing it to a number).
$ws.Range("D3").Value = "'9865486245"

# Move the sheet's active cell/selection to D4.
$ws.Range("D4").Select() | Out-Null
